$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 75) down into the new row 76
# before writing values, so the new row picks up the same date /
# number formatting used by the rest of the table.
$ws.Range("A75:J75").Copy()
$ws.Range("A76:J76").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New data row for 2020-05-25 (serial date 43976)
$ws.Range("A76").Value2 = 43976
$ws.Range("B76").Value2 = 75770
$ws.Range("C76").Value2 = 754
$ws.Range("D76").Value2 = 1469
$ws.Range("E76").Value2 = 0
$ws.Range("F76").Value2 = 9
$ws.Range("G76").Value2 = 2
$ws.Range("H76").Value2 = 6
$ws.Range("I76").Value2 = 108
$ws.Range("J76").Value2 = 1

# Grow the worksheet table / autofilter range to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J76"))

# Match the selection left behind by the editor
$ws.Range("A76:J76").Select()
